$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handback report now reflects that the zh-cn/de-de targets for the
# "650f1eed..." source file are no longer in sync with en-US, so the
# Status text changes everywhere it is shown (Overview summary columns +
# the per-language "Status" column).
$wsOverview.Range("E2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: not in sync with en-US"

$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: not in sync with en-US"

$wsDeDe.Range("C2").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: not in sync with en-US"

# New handback was generated for the "650f1eed..." row (row 2) on both
# language sheets, so its "Correspond Handback DateTime" is refreshed.
$wsZhCn.Range("K2").Value = "2016-10-13 13:56:06"
$wsDeDe.Range("K2").Value = "2016-10-13 13:56:23"

# The longer Status text means the "Status" columns need to be re-fit
# (Overview's zh-cn/de-de summary columns E & F, and the per-language
# sheets' "Status" column C).
$wsOverview.Range("E1").ColumnWidth = 33.4602203369141
$wsOverview.Range("F1").ColumnWidth = 33.4602203369141
$wsZhCn.Range("C1").ColumnWidth = 33.4602203369141
$wsDeDe.Range("C1").ColumnWidth = 33.4602203369141
